# Append three new "sentence" blocks (54's word list, plus new groups 55, 56, 57)
# to the raw-sentence word-frequency sheet, mirroring the existing layout:
#   - a leading run of 6 "blank placeholder" rows (D index 0-5, red fill s=2,
#     empty styled E/F cells) for groups that use the red/green marker rows,
#   - then 14 (or 20) rows carrying D=index, E=word (shared string), F=1,
#   - a blank separator row between groups,
#   - C column holds the group number on the first row of a block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$RED = 255        # fill FFFF0000 (style index 2 in the original workbook)
$GREEN = 5296274  # fill FF92D050 (style index 4 in the original workbook)

# ---------------------------------------------------------------------
# Block 54 (rows 1136-1155): existing block that picks up E/F word data
# from row 1142 onward, and gets red-fill placeholders added for D1136:F1141.
# ---------------------------------------------------------------------
$block54Words = @('[b''again'']', '[b''my'']', '[b''name'']', '[b''is'']', '[b''david'']', '[b''kent'']', '[b''and'']', '[b''i'']', '[b''hope'']', '[b''you'']', '[b''have'']', '[b''a'']', '[b''great'']', '[b''day'']')

for ($i = 0; $i -le 5; $i++) {
    $row = 1136 + $i
    $ws.Cells.Item($row, 4).Value = $i
    $ws.Cells.Item($row, 4).Interior.Color = $RED
    $ws.Cells.Item($row, 5).Interior.Color = $RED
    $ws.Cells.Item($row, 6).Interior.Color = $RED
}

for ($i = 0; $i -lt $block54Words.Length; $i++) {
    $row = 1142 + $i
    $ws.Cells.Item($row, 4).Value = $i + 6
    $ws.Cells.Item($row, 5).Value = $block54Words[$i]
    $ws.Cells.Item($row, 6).Value = 1
}

# ---------------------------------------------------------------------
# Block 55 (rows 1157-1176): brand-new block, C1157 = 55 (green marker row),
# D goes 0-19 with a word + count on every row.
# ---------------------------------------------------------------------
$block55Words = @('[b''encourage'']', '[b''you'']', '[b''to'']', '[b''place'']', '[b''this'']', '[b''display'']', '[b''either'']', '[b''at'']', '[b''your'']', '[b''check'']', '[b''in'']', '[b''or'']', '[b''checkout'']', '[b''counter'']', '[b''or'']', '[b''right'']', '[b''inside'']', '[b''the'']', '[b''treatment'']', '[b''room'']')

$ws.Cells.Item(1157, 3).Value = 55
$ws.Cells.Item(1157, 3).Interior.Color = $GREEN

for ($i = 0; $i -lt $block55Words.Length; $i++) {
    $row = 1157 + $i
    $ws.Cells.Item($row, 4).Value = $i
    $ws.Cells.Item($row, 5).Value = $block55Words[$i]
    $ws.Cells.Item($row, 6).Value = 1
}

# ---------------------------------------------------------------------
# Block 56 (rows 1178-1197): repeats block 54's word list, with the
# same leading red-fill placeholder rows as block 54 (but C1178 itself
# is NOT marked with the green fill).
# ---------------------------------------------------------------------
$block56Words = @('[b''again'']', '[b''my'']', '[b''name'']', '[b''is'']', '[b''david'']', '[b''kent'']', '[b''and'']', '[b''i'']', '[b''hope'']', '[b''you'']', '[b''have'']', '[b''a'']', '[b''great'']', '[b''day'']')

$ws.Cells.Item(1178, 3).Value = 56

for ($i = 0; $i -le 5; $i++) {
    $row = 1178 + $i
    $ws.Cells.Item($row, 4).Value = $i
    $ws.Cells.Item($row, 4).Interior.Color = $RED
    $ws.Cells.Item($row, 5).Interior.Color = $RED
    $ws.Cells.Item($row, 6).Interior.Color = $RED
}

for ($i = 0; $i -lt $block56Words.Length; $i++) {
    $row = 1184 + $i
    $ws.Cells.Item($row, 4).Value = $i + 6
    $ws.Cells.Item($row, 5).Value = $block56Words[$i]
    $ws.Cells.Item($row, 6).Value = 1
}

# ---------------------------------------------------------------------
# Block 57 (rows 1199-1218): brand-new block, D-only (no E/F data yet),
# D goes 0-19, matching the plain "index" pattern used for brand-new,
# not-yet-annotated sentences elsewhere in the sheet.
# ---------------------------------------------------------------------
$ws.Cells.Item(1199, 3).Value = 57

for ($i = 0; $i -le 19; $i++) {
    $row = 1199 + $i
    $ws.Cells.Item($row, 4).Value = $i
}

# ---------------------------------------------------------------------
# Final view state: scroll/select to match where editing ended up.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1186
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F1217").Select()
